$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '28.426.73'
$ws.Range('E2').Value = '  +4.34%  '

Set-TextValue $ws 'D3' '1.794.58'
$ws.Range('E3').Value = '  +1.17%  '

Set-TextValue $ws 'D4' '1.001'
$ws.Range('E4').Value = '  +0.01%  '

Set-TextValue $ws 'D5' '314.57'
$ws.Range('E5').Value = '  +0.51%  '

Set-TextValue $ws 'D6' '1.001'
$ws.Range('E6').Value = '  +0.05%  '

Set-TextValue $ws 'D7' '0.5433'
$ws.Range('E7').Value = '  +3.34%  '

Set-TextValue $ws 'D8' '0.3835'
$ws.Range('E8').Value = '  +4.12%  '

Set-TextValue $ws 'D9' '0.07579'
$ws.Range('E9').Value = '  +3.01%  '

Set-TextValue $ws 'D10' '42.43'
$ws.Range('E10').Value = '  -0.74%  '

Set-TextValue $ws 'D11' '1.123'
$ws.Range('E11').Value = '  +3.22%  '

Set-TextValue $ws 'D12' '1.001'
$ws.Range('E12').Value = '  -0.02%  '

Set-TextValue $ws 'D13' '21.11'
$ws.Range('E13').Value = '  +3.26%  '

Set-TextValue $ws 'D14' '6.189'
$ws.Range('E14').Value = '  +2.09%  '

Set-TextValue $ws 'D15' '7.404'
$ws.Range('E15').Value = '  +6.90%  '

Set-TextValue $ws 'D16' '1.796.11'
$ws.Range('E16').Value = '  +1.58%  '

Set-TextValue $ws 'D17' '91.70'
$ws.Range('E17').Value = '  +3.36%  '

$ws.Range('E18').Value = '  +2.32%  '

Set-TextValue $ws 'D19' '0.06459'
$ws.Range('E19').Value = '  +0.31%  '

Set-TextValue $ws 'D20' '1.001'
$ws.Range('E20').Value = '  +0.04%  '

Set-TextValue $ws 'D21' '17.34'
$ws.Range('E21').Value = '  +3.76%  '

Set-TextValue $ws 'D22' '5.963'
$ws.Range('E22').Value = '  +2.76%  '

Set-TextValue $ws 'D23' '28.407.01'
$ws.Range('E23').Value = '  +4.08%  '

Set-TextValue $ws 'D24' '11.37'
$ws.Range('E24').Value = '  +0.85%  '

Set-TextValue $ws 'D25' '2.121'
$ws.Range('E25').Value = '  +0.55%  '

Set-TextValue $ws 'D26' '159.74'
$ws.Range('E26').Value = '  +2.86%  '

Set-TextValue $ws 'D27' '20.69'

Set-TextValue $ws 'D28' '2.398'
$ws.Range('E28').Value = '  +3.14%  '

Set-TextValue $ws 'D29' '2.002.93'
$ws.Range('E29').Value = '  +1.50%  '

Set-TextValue $ws 'D30' '123.13'
$ws.Range('E30').Value = '  +1.69%  '

Set-TextValue $ws 'D31' '1.121'
$ws.Range('E31').Value = '  +5.73%  '

Set-TextValue $ws 'D32' '0.1022'
$ws.Range('E32').Value = '  +4.60%  '

Set-TextValue $ws 'D33' '5.740'
$ws.Range('E33').Value = '  +3.39%  '

Set-TextValue $ws 'D34' '3.698'
$ws.Range('E34').Value = '  +2.26%  '

Set-TextValue $ws 'D35' '0.2323'
$ws.Range('E35').Value = '  +14.89%  '

Set-TextValue $ws 'D36' '0.06401'
$ws.Range('E36').Value = '  +7.25%  '

Set-TextValue $ws 'D37' '0.02321'
$ws.Range('E37').Value = '  +3.70%  '

$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws 'D38' '8.792'
$ws.Range('E38').Value = '  +8.78%  '

$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws 'D39' '5.145'
$ws.Range('E39').Value = '  +6.43%  '

Set-TextValue $ws 'D40' '11.62'
$ws.Range('E40').Value = '  +3.73%  '

Set-TextValue $ws 'D41' '0.6394'
$ws.Range('E41').Value = '  +4.31%  '

$ws.Range('E42').Value = '  +0.12%  '

Set-TextValue $ws 'D43' '1.159'
$ws.Range('E43').Value = '  +1.76%  '

Set-TextValue $ws 'D44' '1.390'
$ws.Range('E44').Value = '  -2.71%  '

Set-TextValue $ws 'D45' '13.55'
$ws.Range('E45').Value = '  +3.41%  '

Set-TextValue $ws 'D46' '0.5964'
$ws.Range('E46').Value = '  +3.65%  '

Set-TextValue $ws 'D47' '3.675'
$ws.Range('E47').Value = '  +1.46%  '

Set-TextValue $ws 'D48' '126.10'
$ws.Range('E48').Value = '  +4.01%  '

Set-TextValue $ws 'D49' '1.985'
$ws.Range('E49').Value = '  +5.75%  '

Set-TextValue $ws 'D50' '1.150'
$ws.Range('E50').Value = '  +3.19%  '

Set-TextValue $ws 'D51' '0.06888'
$ws.Range('E51').Value = '  +2.92%  '
